$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1512.1428
$ws.Range("I19").Value = 530
$ws.Range("J19").Value = 2248.75
$ws.Range("K19").Value = 530
$ws.Range("L19").Value = 2248.75
$ws.Range("M19").Value = -355
$ws.Range("N19").Value = -2598.75
$ws.Range("H28").Value = 596.1667
$ws.Range("I28").Value = 596.1667
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 596.1667
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -111.1667
$ws.Range("N28").ClearContents()
$ws.Range("H38").Value = 56
$ws.Range("I38").Value = 56
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 168
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 204
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 565
$ws.Range("I58").Value = 565
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1695
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1545
$ws.Range("N58").ClearContents()
$ws.Range("H76").Value = 2737.5
$ws.Range("I76").Value = 2737.5
$ws.Range("K76").Value = 2737.5
$ws.Range("M76").Value = -2422.5
$ws.Range("H79").Value = 2737.5
$ws.Range("I79").Value = 2737.5
$ws.Range("K79").Value = 2737.5
$ws.Range("M79").Value = -1645.5
$ws.Range("H86").Value = 4719.5557
$ws.Range("I86").Value = 3671
$ws.Range("J86").Value = 8389.5
$ws.Range("K86").Value = 3671
$ws.Range("L86").Value = 8389.5
$ws.Range("M86").Value = -2548
$ws.Range("N86").Value = -10635.5
$ws.Range("H89").Value = 4719.5557
$ws.Range("I89").Value = 3671
$ws.Range("J89").Value = 8389.5
$ws.Range("K89").Value = 18355
$ws.Range("L89").Value = 41947.5
$ws.Range("M89").Value = -12739
$ws.Range("N89").Value = -53179.5
$ws.Range("H98").Value = 2387.2856
$ws.Range("I98").Value = 2450.8333
$ws.Range("J98").Value = 2006
$ws.Range("K98").Value = 2450.8333
$ws.Range("L98").Value = 2006
$ws.Range("M98").Value = -952.8332999999998
$ws.Range("N98").Value = -5002
$ws.Range("H122").Value = 2387.2856
$ws.Range("I122").Value = 2450.8333
$ws.Range("J122").Value = 2006
$ws.Range("K122").Value = 7352.499899999999
$ws.Range("L122").Value = 6018
$ws.Range("M122").Value = -4902.499899999999
$ws.Range("N122").Value = -10918
$ws.Range("H132").Value = 772.3
$ws.Range("I132").Value = 772.3
$ws.Range("K132").Value = 2316.9
$ws.Range("M132").Value = 213.1000000000004
$ws.Range("H138").Value = 5404.4614
$ws.Range("I138").Value = 4125.8
$ws.Range("J138").Value = 9666.666999999999
$ws.Range("K138").Value = 12377.4
$ws.Range("L138").Value = 29000.001
$ws.Range("M138").Value = -7237.400000000001
$ws.Range("N138").Value = -39280.001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4083
$ws.Range("I32").Value = 3627.6843
$ws.Range("K32").Value = 3627.6843
$ws.Range("M32").Value = -3340.6843
$ws.Range("H132").Value = 4265.6665
$ws.Range("I132").Value = 4265.6665
$ws.Range("K132").Value = 12796.9995
$ws.Range("M132").Value = -10266.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -12465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 610
$ws.Range("I22").Value = 610
$ws.Range("K22").Value = 610
$ws.Range("M22").Value = -260
$ws.Range("H31").Value = 1230.7097
$ws.Range("I31").Value = 876.3913
$ws.Range("K31").Value = 876.3913
$ws.Range("M31").Value = -581.3913
$ws.Range("H34").Value = 1230.7097
$ws.Range("I34").Value = 876.3913
$ws.Range("K34").Value = 876.3913
$ws.Range("M34").Value = -674.3913
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2230.4211
$ws.Range("J132").Value = 1899.6666
$ws.Range("L132").Value = 5698.9998
$ws.Range("N132").Value = -10758.9998
$ws.Range("H134").Value = 1417.3846
$ws.Range("I134").Value = 1328.3636
$ws.Range("K134").Value = 3985.0908
$ws.Range("M134").Value = -1450.0908

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44024.28
$ws.Range("J2").Value = 84641
$ws.Range("L2").Value = 507846
$ws.Range("N2").Value = -508072
$ws.Range("H23").Value = 669.9
$ws.Range("J23").Value = 399.875
$ws.Range("L23").Value = 1199.625
$ws.Range("N23").Value = -1669.625
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2324
$ws.Range("N94").ClearContents()
$ws.Range("H109").Value = 3181.818
$ws.Range("J109").Value = 3181.818
$ws.Range("L109").Value = 9545.454000000002
$ws.Range("N109").Value = -11625.454
$ws.Range("H111").Value = 499
$ws.Range("I111").Value = 499
$ws.Range("K111").Value = 1497
$ws.Range("M111").Value = 1570
$ws.Range("H121").Value = 3333.3333
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 3333.3333
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 9999.999899999999
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -12619.9999
$ws.Range("H131").Value = 775
$ws.Range("I131").Value = 550
$ws.Range("J131").Value = 1000
$ws.Range("K131").Value = 1650
$ws.Range("L131").Value = 3000
$ws.Range("M131").Value = 3390
$ws.Range("N131").Value = -13080
$ws.Range("H139").Value = 4230.8
$ws.Range("I139").Value = 4230.8
$ws.Range("K139").Value = 12692.4
$ws.Range("M139").Value = -7552.400000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 273.82352
$ws.Range("J2").Value = 472.77777
$ws.Range("L2").Value = 472.77777
$ws.Range("N2").Value = -698.7777699999999
$ws.Range("H70").Value = 2950
$ws.Range("I70").Value = 2950
$ws.Range("K70").Value = 2950
$ws.Range("M70").Value = -2680
$ws.Range("H73").Value = 2950
$ws.Range("I73").Value = 2950
$ws.Range("K73").Value = 2950
$ws.Range("M73").Value = -2014

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 3675
$ws.Range("I81").Value = 3675
$ws.Range("K81").Value = 7350
$ws.Range("M81").Value = -6289
$ws.Range("H84").Value = 3675
$ws.Range("I84").Value = 3675
$ws.Range("K84").Value = 36750
$ws.Range("M84").Value = -31446
$ws.Range("H92").Value = 28900
$ws.Range("J92").Value = 28900
$ws.Range("L92").Value = 28900
$ws.Range("N92").Value = -33892
